$p = $ppt.ActivePresentation
$newDesign = $p.Designs.Add()
$newMaster = $newDesign.SlideMaster
$newLayout = $newMaster.CustomLayouts.Add(1, 1)
$newLayout.Name = "Basic"
# reassign slide to new layout
$p.Slides.Item(1).CustomLayout = $newLayout
Write-Output ("count before delete=" + $p.Designs.Count)
$oldDesign = $p.Designs.Item(1)
$oldDesign.Delete()
Write-Output ("count after delete=" + $p.Designs.Count)
